$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# We apply edits from the END of the document towards the START so that
# paragraph-insert/delete operations never invalidate the 1-based indices
# used by edits that still need to run.
# ---------------------------------------------------------------------------

# --- Hunk 5: simple text substitution (index independent) -----------------
$d.Content.Find.Execute(
    "hər hansısa hack hücumları(XSRF/CSRF attacks)", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "hər hansısa kiber hücumlar(XSRF/CSRF attacks)", 2) | Out-Null

# --- Hunk 4: remove 3 of the 12 blank paragraphs right before "Migrations" -
$d.Paragraphs.Item(123).Range.Delete()
$d.Paragraphs.Item(122).Range.Delete()
$d.Paragraphs.Item(121).Range.Delete()

# --- Hunk 3: split the ConnectionStrings paragraph (Item 110) into 5 ------
$pConn = $d.Paragraphs.Item(110)
$pConn.Range.Text = "“ConnectionStrings” :{"

$pConn.Range.InsertParagraphAfter()
$pE1 = $d.Paragraphs.Item(111)
$pE1.Range.InsertParagraphAfter()

$pDefault = $d.Paragraphs.Item(112)
$pDefault.Range.Text = "“Default” : “Server:ServerName;Database=DatabaseName;Trustesd_Connection=True;”"

$pDefault.Range.InsertParagraphAfter()
$pE2 = $d.Paragraphs.Item(113)
$pE2.Range.InsertParagraphAfter()

$pClose = $d.Paragraphs.Item(114)
$pClose.Range.Text = "}"

# --- Hunk 2: remove 2 of the 6 blank paragraphs right before "View Model" -
$d.Paragraphs.Item(25).Range.Delete()
$d.Paragraphs.Item(24).Range.Delete()

# --- Hunk 1: rewrite the two paragraphs under "Razor View - View Start" ---
$pIntro = $d.Paragraphs.Item(14)
$pIntro.Range.Find.Execute(
    "2 cür olur", $false, $false, $false, $false, $false, $true, 1, $false,
    "2 yerə bölünür", 2) | Out-Null

$pDetail = $d.Paragraphs.Item(16)
$pDetail.Range.Text = "Use Layout’u clickləyib lakin hər hansısa ünvan göstərməyib input’u boş buraxsaq o zaman Razor View Start file’ın olduğu anlaşılır. Razor View Start file View folder daxilində yerləşir. Burada Layout =  _Layout hazır qeyd edilir və bununla da biz Razor View yaradan zaman dediyimiz kimi Use Layout’u clickləyib, Layout’un ünvanını qeyd etməsək belə Layout işə salınır."

$pDetail.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(17)
$pNew1.Range.InsertParagraphAfter()

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
